# Weekly update: a new observation is inserted as row 86 (pushing the
# existing rows 86-218 down to 87-219), which is how "Fruta / hortaliza,
# semanal" commits land new records in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 86; everything below (old rows 86-218)
# shifts down to 87-219, carrying its formatting/styles with it.
$ws.Rows("86:86").Insert()

# Populate the newly inserted row with the new observation's data.
$ws.Cells.Item(86, 1).Value = 7
$ws.Cells.Item(86, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(86, 3).Value = "Ñuble"
$ws.Cells.Item(86, 4).Value = 44540
$ws.Cells.Item(86, 5).Value = 16
$ws.Cells.Item(86, 6).Value = 100114013
$ws.Cells.Item(86, 7).Value = "Zanahoria"
$ws.Cells.Item(86, 8).Value = "Sin especificar"
$ws.Cells.Item(86, 9).Value = "Primera"
$ws.Cells.Item(86, 10).Value = 100
$ws.Cells.Item(86, 11).Value = 8000
$ws.Cells.Item(86, 12).Value = 8500
$ws.Cells.Item(86, 13).Value = 8250
$ws.Cells.Item(86, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(86, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(86, 16).Value = 412
$ws.Cells.Item(86, 17).Value = 20
$ws.Cells.Item(86, 18).Value = "Hortaliza"

# Make sure the date column keeps the same date/time number format as the
# rest of column D (the Insert already copies it, but set it explicitly
# to be safe).
$ws.Cells.Item(86, 4).NumberFormat = $ws.Cells.Item(87, 4).NumberFormat
